# Loan RBI, Variable Instalments
# ---------------------------------------------------------------
# The "Repayment schedule" sheet gets a new (blank-header) column
# inserted before the existing "Late" column (old column N), which
# pushes "Late" / "heading" (Date) / "Disbursement" one column to
# the right (N->O, O->P, P->Q). The sheet also becomes the active
# tab/selection in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make this the active sheet (it becomes the selected tab on save).
$ws.Activate()

# Insert a new blank column at N - everything from N onward (Late,
# the trailing Date heading, Disbursement) shifts right by one.
$ws.Range("N1").EntireColumn.Insert()

# The newly inserted column picks up the width of its neighbour
# (11 characters, matching column M) instead of Excel's generic
# default width.
$ws.Columns.Item(14).ColumnWidth = 10.166666666666666

# Update the active selection to match the post-edit cursor position.
$ws.Range("R8").Select() | Out-Null
